# MASSACHUSETTS_2016.xlsx cleanup
# - rename header columns to short snake_case names
# - title-case the Spanish state/municipality names in columns A and B
# - correct a handful of 1-ULP floating point roundings in column D
# - drop the trailing footnote rows (381-480), shrinking the used range to A1:D379

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header row rename
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case every non-empty text cell in columns A (state) and B (municipality)
#    for the data rows (2-379). Matches .title()-style casing: "de" -> "De",
#    "del" -> "Del", "el" -> "El", "los" -> "Los", "y" -> "Y", etc.
$ti = (Get-Culture).TextInfo
for ($r = 2; $r -le 379; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value()
    if ($valA -ne $null -and $valA -ne "") {
        $cellA.Value = $ti.ToTitleCase($valA)
    }
    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value()
    if ($valB -ne $null -and $valB -ne "") {
        $cellB.Value = $ti.ToTitleCase($valB)
    }
}

# 3) A few D-column percentages were re-rounded by 1 ULP on regeneration.
$ws.Range("D21").Value = 0.009174311926605503
$ws.Range("D56").Value = 0.009174311926605503
$ws.Range("D74").Value = 0.009174311926605503
$ws.Range("D78").Value = 0.009174311926605503
$ws.Range("D81").Value = 0.09403669724770644
$ws.Range("D180").Value = 0.09403669724770644
$ws.Range("D202").Value = 0.009174311926605503
$ws.Range("D252").Value = 0.009174311926605503

# 4) Drop the old footnote rows that used to live below the data (381-385 and
#    476-480), shrinking the sheet's used range down to A1:D379.
$ws.Rows("380:480").Delete()
